$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.344.80'
$ws.Range("E2").Value = '  -0.59%  '

$ws.Range("D3").Value = '2.390.67'
$ws.Range("E3").Value = '  -3.83%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''549.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("D6").Value = '''141.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.11%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '''0.535'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -11.03%  '

$ws.Range("D9").Value = '2.389.76'
$ws.Range("E9").Value = '  -3.86%  '

$ws.Range("D10").Value = '''0.106'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.66%  '

$ws.Range("E11").Value = '  +0.08%  '

$ws.Range("D12").Value = '''5.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.23%  '

$ws.Range("D13").Value = '''0.348'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.07%  '

$ws.Range("D14").Value = '''25.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.45%  '

$ws.Range("D15").Value = '2.821.02'
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("E16").Value = '  -2.41%  '

$ws.Range("D17").Value = '61.108.04'
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").Value = '2.389.88'
$ws.Range("E18").Value = '  -3.57%  '

$ws.Range("D19").Value = '''10.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.51%  '

$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("D21").Value = '''318.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '

$ws.Range("D22").Value = '''6.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.16%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").Value = '''1.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '

$ws.Range("D25").Value = '''63.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.67%  '

$ws.Range("D26").Value = '''8.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.62%  '

$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").Value = '2.507.16'
$ws.Range("E28").Value = '  -3.56%  '

$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '''530.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.35%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0927'
$ws.Range("E30").Value = '  -8.64%  '

$ws.Range("D31").Value = '''8.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("E32").Value = '  -5.83%  '

$ws.Range("E33").Value = '  -4.04%  '

$ws.Range("E34").Value = '  -4.07%  '

$ws.Range("E35").Value = '  -0.86%  '

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("E37").Value = '  -7.51%  '

$ws.Range("E38").Value = '  -4.57%  '

$ws.Range("D39").Value = '''0.375'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.66%  '

$ws.Range("D40").Value = '''1.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.21%  '

$ws.Range("D41").Value = '''18.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.93%  '

$ws.Range("D42").Value = '''140.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.85%  '

$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("D44").Value = '''40.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '''3.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.15%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''140.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.32%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.65%  '

$ws.Range("D48").Value = '''20.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.64%  '

$ws.Range("D49").Value = '''0.0520'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.39%  '

$ws.Range("D50").Value = '''0.577'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.74%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.0908'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.04%  '
